$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 is new: seed its formatting (date style, bold headers border, etc.)
# by copying the format pattern from the row directly above it (row 18),
# matching the existing per-row style already used for column A (date column).
$ws.Range("A18").Copy($ws.Range("A19"))

# Update dimension data: each row now holds (date, y0 year, y0_forecast, y1 year, y1_forecast)
$ws.Cells.Item(2, 1).Value = 39400
$ws.Cells.Item(2, 2).Value = 2007
$ws.Cells.Item(2, 3).Value = 4.930115226412335
$ws.Cells.Item(2, 4).Value = 2008
$ws.Cells.Item(2, 5).Value = 1.402909115503936

$ws.Cells.Item(3, 1).Value = 39765
$ws.Cells.Item(3, 2).Value = 2008
$ws.Cells.Item(3, 3).Value = 1.457587285166628
$ws.Cells.Item(3, 4).Value = 2009
$ws.Cells.Item(3, 5).Value = 2.807231216534278

$ws.Cells.Item(4, 1).Value = 40130
$ws.Cells.Item(4, 2).Value = 2009
$ws.Cells.Item(4, 3).Value = -0.9140166223623569
$ws.Cells.Item(4, 4).Value = 2010
$ws.Cells.Item(4, 5).Value = 1.821983295885121

$ws.Cells.Item(5, 1).Value = 40494
$ws.Cells.Item(5, 2).Value = 2010
$ws.Cells.Item(5, 3).Value = 2.585942866987878
$ws.Cells.Item(5, 4).Value = 2011
$ws.Cells.Item(5, 5).Value = 2.722861752007866

$ws.Cells.Item(6, 1).Value = 40862
$ws.Cells.Item(6, 2).Value = 2011
$ws.Cells.Item(6, 3).Value = 4.253963781362402
$ws.Cells.Item(6, 4).Value = 2012
$ws.Cells.Item(6, 5).Value = 2.878414118480799

$ws.Cells.Item(7, 1).Value = 41228
$ws.Cells.Item(7, 2).Value = 2012
$ws.Cells.Item(7, 3).Value = 1.752870900283909
$ws.Cells.Item(7, 4).Value = 2013
$ws.Cells.Item(7, 5).Value = 3.144721336271927

$ws.Cells.Item(8, 1).Value = 41592
$ws.Cells.Item(8, 2).Value = 2013
$ws.Cells.Item(8, 3).Value = -1.479696720105139
$ws.Cells.Item(8, 4).Value = 2014
$ws.Cells.Item(8, 5).Value = 2.238623952069552

$ws.Cells.Item(9, 1).Value = 41957
$ws.Cells.Item(9, 2).Value = 2014
$ws.Cells.Item(9, 3).Value = 3.900127535411246
$ws.Cells.Item(9, 4).Value = 2015
$ws.Cells.Item(9, 5).Value = -0.6155071485167585

$ws.Cells.Item(10, 1).Value = 42321
$ws.Cells.Item(10, 2).Value = 2015
$ws.Cells.Item(10, 3).Value = 0.03947433952959933
$ws.Cells.Item(10, 4).Value = 2016
$ws.Cells.Item(10, 5).Value = 1.182212550358064

$ws.Cells.Item(11, 1).Value = 42689
$ws.Cells.Item(11, 2).Value = 2016
$ws.Cells.Item(11, 3).Value = 2.192778679161944
$ws.Cells.Item(11, 4).Value = 2017
$ws.Cells.Item(11, 5).Value = -0.5835597102573198

$ws.Cells.Item(12, 1).Value = 43053
$ws.Cells.Item(12, 2).Value = 2017
$ws.Cells.Item(12, 3).Value = 3.40836448860673
$ws.Cells.Item(12, 4).Value = 2018
$ws.Cells.Item(12, 5).Value = 2.34299484087257

$ws.Cells.Item(13, 1).Value = 43418
$ws.Cells.Item(13, 2).Value = 2018
$ws.Cells.Item(13, 3).Value = 2.799070570134488
$ws.Cells.Item(13, 4).Value = 2019
$ws.Cells.Item(13, 5).Value = 3.056075254339996

$ws.Cells.Item(14, 1).Value = 43783
$ws.Cells.Item(14, 2).Value = 2019
$ws.Cells.Item(14, 3).Value = 4.195393191694419
$ws.Cells.Item(14, 4).Value = 2020
$ws.Cells.Item(14, 5).Value = 2.652948310315506

$ws.Cells.Item(15, 1).Value = 44159
$ws.Cells.Item(15, 2).Value = 2020
$ws.Cells.Item(15, 3).Value = 1.666553973046048
$ws.Cells.Item(15, 4).Value = 2021
$ws.Cells.Item(15, 5).Value = -0.9999522486825452

$ws.Cells.Item(16, 1).Value = 44525
$ws.Cells.Item(16, 2).Value = 2021
$ws.Cells.Item(16, 3).Value = 1.879266440112803
$ws.Cells.Item(16, 4).Value = 2022
$ws.Cells.Item(16, 5).Value = -0.1343977949472275

$ws.Cells.Item(17, 1).Value = 44890
$ws.Cells.Item(17, 2).Value = 2022
$ws.Cells.Item(17, 3).Value = -2.620683231370946
$ws.Cells.Item(17, 4).Value = 2023
$ws.Cells.Item(17, 5).Value = -2.83913279674276

$ws.Cells.Item(18, 1).Value = 45254
$ws.Cells.Item(18, 2).Value = 2023
$ws.Cells.Item(18, 3).Value = -3.036556262700274
$ws.Cells.Item(18, 4).Value = 2024
$ws.Cells.Item(18, 5).Value = -1.803491225663911

$ws.Cells.Item(19, 1).Value = 45618
$ws.Cells.Item(19, 2).Value = 2024
$ws.Cells.Item(19, 3).Value = -2.953443685011514
$ws.Cells.Item(19, 4).Value = 2025
$ws.Cells.Item(19, 5).Value = -2.75492543068685
